# The "OUTLINE" slide has a bullet that currently reads:
#   "Making sense of the data collection (Live Session)"
# Split it into two runs so the leading clause becomes
#   "Integrating data from multiple sources ("
# while the trailing "Live Session)" remains as its own run, unchanged.

$targetParaText = "Making sense of the data collection (Live Session)"
$oldLead        = "Making sense of the data collection ("
$newLead        = "Integrating data from multiple sources ("

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }

        $tf = $shape.TextFrame
        if (-not $tf.HasText) { continue }

        $tr = $tf.TextRange
        $paraCount = $tr.Paragraphs().Count

        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi)
            $paraText = $para.Text.TrimEnd("`r")

            if ($paraText -eq $targetParaText) {
                # Replace only the leading portion; the remaining
                # "Live Session)" text/run is left exactly as-is.
                $leadRange = $tr.Characters($para.Start, $oldLead.Length)
                $leadRange.Text = $newLead
            }
        }
    }
}
